$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 12:50"
$ws.Range("A6").Value = "Bizkaia/Vizcaya"
$ws.Range("B6").Value = 3102
$ws.Range("C6").Value = 1796
$ws.Range("D6").Value = 2154
$ws.Range("E6").Value = 148

$ws.Range("A7").Value = "Valencia/Valencia"
$ws.Range("B7").Value = 2874
$ws.Range("C7").Value = 138
$ws.Range("D7").Value = 2578
$ws.Range("E7").Value = 158

$ws.Range("A8").Value = "Castilla-La Mancha"
$ws.Range("B8").Value = 2780
$ws.Range("C8").Value = 71
$ws.Range("D8").Value = 2446
$ws.Range("E8").Value = 263

$ws.Range("A10").Value = "Araba/Alava"
$ws.Range("B10").Value = 2105
$ws.Range("C10").Value = 1796
$ws.Range("D10").Value = 1372
$ws.Range("E10").Value = 133

$ws.Range("A12").Value = "Alacant/Alicante"
$ws.Range("B12").Value = 2021
$ws.Range("C12").Value = 53
$ws.Range("D12").Value = 1821
$ws.Range("E12").Value = 147

$ws.Range("A14").Value = "Zaragoza"
$ws.Range("B14").Value = 1792
$ws.Range("C14").Value = 165
$ws.Range("D14").Value = 1520
$ws.Range("E14").Value = 107

$ws.Range("A15").Value = "A Coruña"
$ws.Range("B15").Value = 1687
$ws.Range("C15").Value = 187
$ws.Range("D15").Value = 1573
$ws.Range("E15").Value = 54

$ws.Range("A16").Value = "Albacete"
$ws.Range("B16").Value = 1537
$ws.Range("C16").Value = 252
$ws.Range("D16").Value = 1339
$ws.Range("E16").Value = 133

$ws.Range("A17").Value = "Toledo"
$ws.Range("B17").Value = 1426
$ws.Range("C17").Value = 252
$ws.Range("D17").Value = 1207
$ws.Range("E17").Value = 155

$ws.Range("A18").Value = "Pontevedra"
$ws.Range("B18").Value = 1380
$ws.Range("C18").Value = 187
$ws.Range("D18").Value = 1294
$ws.Range("E18").Value = 21

$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 1321
$ws.Range("C19").Value = 83
$ws.Range("D19").Value = 1169
$ws.Range("E19").Value = 69

$ws.Range("A20").Value = "Gran Canaria"
$ws.Range("B20").Value = 1262
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 320
$ws.Range("E20").Value = 11

$ws.Range("A21").Value = "La Palma"
$ws.Range("B21").Value = 1262
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = 2

$ws.Range("A22").Value = "Lanzarote"
$ws.Range("B22").Value = 1262
$ws.Range("C22").Value = 32
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = 3

$ws.Range("A23").Value = "Fuerteventura"
$ws.Range("B23").Value = 1262
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = 0

$ws.Range("A24").Value = "La Gomera"
$ws.Range("B24").Value = 1262
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 0

$ws.Range("A25").Value = "El Hierro"
$ws.Range("B25").Value = 1262
$ws.Range("C25").Value = 32
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0

$ws.Range("A26").Value = "Asturias"
$ws.Range("B26").Value = 1236
$ws.Range("C26").Value = 90
$ws.Range("D26").Value = 1091
$ws.Range("E26").Value = 55

$ws.Range("A27").Value = "Cantabria"
$ws.Range("B27").Value = 1171
$ws.Range("C27").Value = 35
$ws.Range("D27").Value = 1099
$ws.Range("E27").Value = 37

$ws.Range("A28").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B28").Value = 1113
$ws.Range("C28").Value = 1796
$ws.Range("D28").Value = 673
$ws.Range("E28").Value = 44

$ws.Range("A29").Value = "Salamanca"
$ws.Range("B29").Value = 1078
$ws.Range("C29").Value = 181
$ws.Range("D29").Value = 787
$ws.Range("E29").Value = 110

$ws.Range("A30").Value = "Caceres"
$ws.Range("B30").Value = 1067
$ws.Range("C30").Value = 29
$ws.Range("D30").Value = 924
$ws.Range("E30").Value = 114

$ws.Range("A31").Value = "Sevilla"
$ws.Range("B31").Value = 1052
$ws.Range("C31").Value = 18
$ws.Range("D31").Value = 1000
$ws.Range("E31").Value = 34

$ws.Range("A39").Value = "Castello/Castellon"
$ws.Range("B39").Value = 613
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 570
$ws.Range("E39").Value = 34

$ws.Range("A40").Value = "Jaen"
$ws.Range("B40").Value = 599
$ws.Range("C40").Value = 17
$ws.Range("D40").Value = 559
$ws.Range("E40").Value = 23

$ws.Range("A43").Value = "Ourense"
$ws.Range("B43").Value = 570
$ws.Range("C43").Value = 187
$ws.Range("D43").Value = 520
$ws.Range("E43").Value = 12

$ws.Range("A44").Value = "Badajoz"
$ws.Range("B44").Value = 561
$ws.Range("C44").Value = 62
$ws.Range("D44").Value = 480
$ws.Range("E44").Value = 19

$ws.Range("A45").Value = "Soria"
$ws.Range("B45").Value = 550
$ws.Range("C45").Value = 71
$ws.Range("D45").Value = 442
$ws.Range("E45").Value = 37

$ws.Range("A46").Value = "Cadiz"
$ws.Range("B46").Value = 507
$ws.Range("C46").Value = 10
$ws.Range("D46").Value = 484
$ws.Range("E46").Value = 13

$ws.Range("A48").Value = "Lugo"
$ws.Range("B48").Value = 402
$ws.Range("C48").Value = 187
$ws.Range("D48").Value = 371
$ws.Range("E48").Value = 7

$ws.Range("A49").Value = "Palencia"
$ws.Range("B49").Value = 293
$ws.Range("C49").Value = 33
$ws.Range("D49").Value = 242
$ws.Range("E49").Value = 18

$ws.Range("A51").Value = "Huesca"
$ws.Range("B51").Value = 244
$ws.Range("C51").Value = 23
$ws.Range("D51").Value = 207
$ws.Range("E51").Value = 14

$ws.Range("A52").Value = "Teruel"
$ws.Range("B52").Value = 236
$ws.Range("C52").Value = 16
$ws.Range("D52").Value = 205
$ws.Range("E52").Value = 15

$ws.Range("A53").Value = "Almeria"
$ws.Range("B53").Value = 223
$ws.Range("C53").Value = 6
$ws.Range("D53").Value = 203
$ws.Range("E53").Value = 14
